# Project DesignFirst is saved. Author: admin. Type: SAVE.
# D10 on the "Rules" sheet is updated from 21 to 100 (numeric), matching
# the sibling cell C10's value/type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = 100
